$wb = $excel.ActiveWorkbook

# --- CreatedSuites: append the newly created suite rows ---
$wsCreated = $wb.Worksheets.Item("CreatedSuites")
$wsCreated.Range("A2").Value = "Suite-22:49:58"
$wsCreated.Range("A3").Value = "Suite-23:40:49"
$wsCreated.Range("A4").Value = "Royal-23:42:05"
$wsCreated.Range("A5").Value = "Royal-23:43:04"

# --- UpdatedSuites: append the newly updated suite rows ---
$wsUpdated = $wb.Worksheets.Item("UpdatedSuites")
$wsUpdated.Range("A2").Value = "Suite-23:40:49"
$wsUpdated.Range("B2").Value = "Royal-23:42:05"
$wsUpdated.Range("A3").Value = "Royal-23:42:05"
$wsUpdated.Range("B3").Value = "Royal-23:43:04"

# --- Switch the active/selected sheet & selection over to CreatedSuites ---
$wsCreated.Activate()
$wsCreated.Range("B7").Select()
